# Amend corrected label annotations
# Normalizes/reorders the "labels" (column F) values for several rows:
# lower-cases category names and re-orders pipe-delimited multi-labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = "ppe"
$ws.Range("F13").Value = "ppe"
$ws.Range("F15").Value = "application instructions || env warning - species || pollinator"
$ws.Range("F16").Value = "application instructions || env warning - species || pollinator"
$ws.Range("F17").Value = "32_physical_and_chemical_hazards"
$ws.Range("F21").Value = "use restrictions"
$ws.Range("F22").Value = "application instructions"
$ws.Range("F24").Value = "application instructions"
$ws.Range("F25").Value = "application instructions"
$ws.Range("F27").Value = "application instructions"
$ws.Range("F28").Value = "safety procedures || application instructions"
$ws.Range("F29").Value = "safety procedures"
$ws.Range("F30").Value = "safety procedures || application instructions"
$ws.Range("F34").Value = "154_pesticide_storage"
$ws.Range("F60").Value = "application instructions || env warning - species || pollinator"
$ws.Range("F61").Value = "32_physical_and_chemical_hazards"
$ws.Range("F65").Value = "use restrictions"
$ws.Range("F66").Value = "application instructions"
$ws.Range("F71").Value = "safety procedures || application instructions"
$ws.Range("F72").Value = "safety procedures || application instructions"
$ws.Range("F76").Value = "154_pesticide_storage"
